# Add a new "total" worksheet (sum per year, mirroring the existing
# "average" worksheet which holds the per-year mean) after the
# "average" sheet.

$wb = $excel.ActiveWorkbook
$wsMain = $wb.Worksheets.Item("main")
$wsAvg = $wb.Worksheets.Item("average")

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "total"

# Reuse the header row's formatting (bold, centered, bordered) from
# the "average" sheet so the existing style is shared instead of
# minting a new one, then fill in the same column labels (these map
# back onto the existing shared strings, so no new unique strings are
# created either).
$wsAvg.Range("B1:G1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)

# Reuse the bold/centered/bordered style used on column A (the row
# index column) on the "average" sheet for the same column here.
$wsAvg.Range("A2:A11").Copy()
$ws.Range("A2:A11").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$ws.Range("B1").Value = "year"
$ws.Range("C1").Value = "Temperature Change C"
$ws.Range("D1").Value = "annual co2 emmisions"
$ws.Range("E1").Value = "GDP `$"
$ws.Range("F1").Value = "Green Bonds Issuance `$"
$ws.Range("G1").Value = "population"

# year -> [temperature change C, annual co2 emissions, GDP $,
#          Green Bonds Issuance $, population] totals (sum across all
# countries for that year)
$yearData = @(
    @(2012, 56.51300000000001, 4284036292, 18899062187231.49, 0.6503946, 591702029),
    @(2013, 39.975, 4187930277, 19710492567527.39, 4.1503489, 593142931),
    @(2014, 79.292, 3958391776, 20293345323920.78, 16.3897677, 594565062),
    @(2015, 68.438, 3960306862, 17839800156187.61, 27.0352689, 595963744),
    @(2016, 72.72999999999998, 3957305821.9, 17923128827861.67, 31.1519053, 597279081),
    @(2017, 60.252, 3960139110, 18847422473370.02, 65.9753953, 598240534),
    @(2018, 81.67700000000002, 3887174521, 20369250013115.85, 82.98431329999998, 599120957),
    @(2019, 79.45500000000001, 3717680253, 20051444797403.53, 151.7125718, 599504550),
    @(2020, 92.11499999999999, 3370241592, 19506692418091.7, 168.0720224, 599695156),
    @(2021, 59.56900000000001, 3536737297, 21977425645774.62, 318.1558474999999, 598984930)
)

$r = 2
$idx = 0
foreach ($row in $yearData) {
    $ws.Cells.Item($r, 1).Value = $idx
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $idx++
    $r++
}

# Restore "main" as the active/selected sheet, matching the original
# workbook's tab selection.
$wsMain.Activate()
